$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had a Name/Phone contact list; replace it with a
# Name/Orders table (Mamun, Nahid, Limon) and drop the old Phone column.
$ws.Range("C1:C2").ClearContents()

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Orders"
$ws.Range("A2").Value = "Mamun"
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = "Nahid"
$ws.Range("B3").Value = 3
$ws.Range("A4").Value = "Limon"
$ws.Range("B4").Value = 3

# Center the new table (horizontally + vertically) and spread that
# formatting across the whole data range via copy/paste-special so it
# only adds a single combined style instead of one style per property.
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").Copy()
$ws.Range("A1:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final on-screen state: zoomed in to 205% with the
# cursor resting on I7.
$excel.ActiveWindow.Zoom = 205
$null = $ws.Range("I7").Select()
